$d = $word.ActiveDocument

# The second paragraph holds a Word field (begin / instrText* / end) that
# encodes the M2Doc query `m:'doc.html'.fromHTMLURI()`. Replace the whole
# paragraph's contents with plain-text runs (splitting the text the same
# way the field code was split across runs), turning the field into
# literal `{ ... }` text while keeping the bookmark in place.
$p = $d.Paragraphs(2)
$r = $p.Range

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
    '<pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
    '<w:body>' +
    '<w:p w:rsidR="00C52979" w:rsidRDefault="00C52979" w:rsidP="00F5495F">' +
    '<w:r><w:t>{</w:t></w:r>' +
    '<w:r><w:t>m</w:t></w:r>' +
    '<w:r><w:t>:</w:t></w:r>' +
    '<w:r><w:t>''</w:t></w:r>' +
    '<w:r><w:t>doc.html</w:t></w:r>' +
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
    '<w:bookmarkEnd w:id="0"/>' +
    '<w:r><w:t>''.fromHTMLURI()</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">}</w:t></w:r>' +
    '</w:p>' +
    '</w:body>' +
    '</w:document>' +
    '</pkg:xmlData>' +
    '</pkg:part>' +
    '</pkg:package>'

[void]$r.InsertXML($xml)
